# -----------------------------------------------------------------------
# epexspot_prices.xlsx automated update: add 04-sep data point to all sheets
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# --- "Prix Spot" sheet: add new column CE ("04-sep") -----------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell CE1: set text, then copy CD1's formatting (bold/border/centered)
# onto it via copy + paste-special so it reuses the existing header style.
$wsPrix.Range("CE1").Value = "04-sep"
$wsPrix.Range("CD1").Copy()
$wsPrix.Range("CE1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for column CE, rows 2-25
$wsPrix.Range("CE2").Value = 19.84
$wsPrix.Range("CE3").Value = 14.81
$wsPrix.Range("CE4").Value = 14.3
$wsPrix.Range("CE5").Value = 5.96
$wsPrix.Range("CE6").Value = 5.99
$wsPrix.Range("CE7").Value = 11.49
$wsPrix.Range("CE8").Value = 42.15
$wsPrix.Range("CE9").Value = 52.04
$wsPrix.Range("CE10").Value = 55.89
$wsPrix.Range("CE11").Value = 49.6
$wsPrix.Range("CE12").Value = 18.07
$wsPrix.Range("CE13").Value = 5.79
$wsPrix.Range("CE14").Value = 3.2
$wsPrix.Range("CE15").Value = 1.72
$wsPrix.Range("CE16").Value = 1.6
$wsPrix.Range("CE17").Value = 1.55
$wsPrix.Range("CE18").Value = 3.7
$wsPrix.Range("CE19").Value = 7.01
$wsPrix.Range("CE20").Value = 28.19
$wsPrix.Range("CE21").Value = 62.19
$wsPrix.Range("CE22").Value = 102.5
$wsPrix.Range("CE23").Value = 80.74
$wsPrix.Range("CE24").Value = 84.8
$wsPrix.Range("CE25").Value = 62.18

# --- "Gaz" sheet: append row 80 (2025-09-02) --------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force the date column to stay plain text (matches existing A2:A79 cells)
# instead of being auto-parsed into a date serial number, then drop the
# temporary text number-format so the cell keeps the sheet default style.
$wsGaz.Range("A80").NumberFormat = "@"
$wsGaz.Range("A80").Value = "2025-09-02"
$wsGaz.Range("A80").ClearFormats()
$wsGaz.Range("B80").Value = 30.5

# --- "CO2" sheet: append row 80 (2025-09-02) ---------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A80").NumberFormat = "@"
$wsCo2.Range("A80").Value = "2025-09-02"
$wsCo2.Range("A80").ClearFormats()
$wsCo2.Range("B80").Value = 73.25
